# "update to manual status column" — the fastq QC sheet's manualStatus
# column (I) held a bare status code (4). Switch it to the bracketed
# text form "[4]" used going forward, for the three rows that had a
# status set (rows 5-7). Those rows also get a touch more breathing
# room / the fastqFileName column (F) is widened to fit its long file
# names, and the active selection moves onto the column that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# manualStatus: numeric 4 -> text "[4]"
$ws.Range("I5").Value = "[4]"
$ws.Range("I6").Value = "[4]"
$ws.Range("I7").Value = "[4]"

# Rows touched by the manualStatus edit get a slightly tighter row height
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Widen the fastqFileName column so the long file names are visible
$ws.Columns.Item(6).ColumnWidth = 50.99

# Leave the selection on the cell that was last edited
$ws.Range("I7").Select()
